# Daily attendance processing - rotate "Recorded By" names in column G.
# For each row, if the "Recorded By" cell contains multiple comma-separated
# names and the first name is NOT "System", move that first name to the end
# of the list (left-rotate by one). Cells that are a single name, or whose
# first name already is "System", are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-RecordedBy($value) {
    if ($null -eq $value) {
        return $value
    }

    $parts = $value.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    if ($trimmed.Length -le 1) {
        return $value
    }

    if ($trimmed[0] -eq "System") {
        return $value
    }

    $first = $trimmed[0]
    $rest = $trimmed[1..($trimmed.Length - 1)]
    $newParts = $rest + @($first)
    return [string]::Join(", ", $newParts)
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($null -eq $current -or $current -eq "") {
        continue
    }
    $updated = Rotate-RecordedBy $current
    if ($updated -ne $current) {
        $cell.Value2 = $updated
    }
}
